$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.005243333333333
$ws.Range("H2").Value = 3.01573
$ws.Range("I2").Value = 0.07224874268505826
$ws.Range("J2").Value = 0.07224874268505825
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 1.378475333333333
$ws.Range("N2").Value = 4.135426
$ws.Range("O2").Value = 0.05609715574531157
$ws.Range("P2").Value = 0.05609715574531156
$ws.Range("Q2").Value = 1.385703138997778
$ws.Range("R2").Value = 12.47132825098
$ws.Range("S2").Value = 0.004052948970806653
$ws.Range("T2").Value = 0.004052948970806652

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.005243333333333
$ws.Range("H3").Value = 3.01573
$ws.Range("I3").Value = 0.07224874268505826
$ws.Range("J3").Value = 0.07224874268505825
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 13.16176133333333
$ws.Range("N3").Value = 39.485284
$ws.Range("O3").Value = 0.5356188518899525
$ws.Range("P3").Value = 0.5356188518899525
$ws.Range("Q3").Value = 13.23077283525778
$ws.Range("R3").Value = 119.07695551732
$ws.Range("S3").Value = 0.03869778860746351
$ws.Range("T3").Value = 0.03869778860746351

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.005243333333333
$ws.Range("H4").Value = 3.01573
$ws.Range("I4").Value = 0.07224874268505826
$ws.Range("J4").Value = 0.07224874268505825
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.829094
$ws.Range("N4").Value = 2.487282
$ws.Range("O4").Value = 0.03374004171190829
$ws.Range("P4").Value = 0.03374004171190828
$ws.Range("Q4").Value = 0.8334412162066666
$ws.Range("R4").Value = 7.500970945859999
$ws.Range("S4").Value = 0.002437675591826794
$ws.Range("T4").Value = 0.002437675591826794

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1.005243333333333
$ws.Range("H5").Value = 3.01573
$ws.Range("I5").Value = 0.07224874268505826
$ws.Range("J5").Value = 0.07224874268505825
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 9.203668
$ws.Range("N5").Value = 27.611004
$ws.Range("O5").Value = 0.3745439506528278
$ws.Range("P5").Value = 0.3745439506528276
$ws.Range("Q5").Value = 9.251925899213333
$ws.Range("R5").Value = 83.26733309292
$ws.Range("S5").Value = 0.02706032951496131
$ws.Range("T5").Value = 0.0270603295149613

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 10.25983933333333
$ws.Range("H6").Value = 30.779518
$ws.Range("I6").Value = 0.7373940889775011
$ws.Range("J6").Value = 0.737394088977501
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 1.378475333333333
$ws.Range("N6").Value = 4.135426
$ws.Range("O6").Value = 0.05609715574531157
$ws.Range("P6").Value = 0.05609715574531156
$ws.Range("Q6").Value = 14.14293544496311
$ws.Range("R6").Value = 127.286419004668
$ws.Range("S6").Value = 0.04136571105504302
$ws.Range("T6").Value = 0.041365711055043

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 10.25983933333333
$ws.Range("H7").Value = 30.779518
$ws.Range("I7").Value = 0.7373940889775011
$ws.Range("J7").Value = 0.737394088977501
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 13.16176133333333
$ws.Range("N7").Value = 39.485284
$ws.Range("O7").Value = 0.5356188518899525
$ws.Range("P7").Value = 0.5356188518899525
$ws.Range("Q7").Value = 135.0375566236791
$ws.Range("R7").Value = 1215.338009613112
$ws.Range("S7").Value = 0.3949621753285666
$ws.Range("T7").Value = 0.3949621753285666

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 10.25983933333333
$ws.Range("H8").Value = 30.779518
$ws.Range("I8").Value = 0.7373940889775011
$ws.Range("J8").Value = 0.737394088977501
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.829094
$ws.Range("N8").Value = 2.487282
$ws.Range("O8").Value = 0.03374004171190829
$ws.Range("P8").Value = 0.03374004171190828
$ws.Range("Q8").Value = 8.506371232230666
$ws.Range("R8").Value = 76.55734109007599
$ws.Range("S8").Value = 0.0248797073202155
$ws.Range("T8").Value = 0.02487970732021549

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 10.25983933333333
$ws.Range("H9").Value = 30.779518
$ws.Range("I9").Value = 0.7373940889775011
$ws.Range("J9").Value = 0.737394088977501
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 9.203668
$ws.Range("N9").Value = 27.611004
$ws.Range("O9").Value = 0.3745439506528278
$ws.Range("P9").Value = 0.3745439506528276
$ws.Range("Q9").Value = 94.42815495734133
$ws.Range("R9").Value = 849.853394616072
$ws.Range("S9").Value = 0.2761864952736761
$ws.Range("T9").Value = 0.2761864952736759

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.383875
$ws.Range("H10").Value = 1.151625
$ws.Range("I10").Value = 0.02758982345723265
$ws.Range("J10").Value = 0.02758982345723265
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 1.378475333333333
$ws.Range("N10").Value = 4.135426
$ws.Range("O10").Value = 0.05609715574531157
$ws.Range("P10").Value = 0.05609715574531156
$ws.Range("Q10").Value = 0.5291622185833333
$ws.Range("R10").Value = 4.76245996725
$ws.Range("S10").Value = 0.001547710623466031
$ws.Range("T10").Value = 0.00154771062346603

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.383875
$ws.Range("H11").Value = 1.151625
$ws.Range("I11").Value = 0.02758982345723265
$ws.Range("J11").Value = 0.02758982345723265
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 13.16176133333333
$ws.Range("N11").Value = 39.485284
$ws.Range("O11").Value = 0.5356188518899525
$ws.Range("P11").Value = 0.5356188518899525
$ws.Range("Q11").Value = 5.052471131833333
$ws.Range("R11").Value = 45.4722401865
$ws.Range("S11").Value = 0.01477762956400944
$ws.Range("T11").Value = 0.01477762956400944

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.383875
$ws.Range("H12").Value = 1.151625
$ws.Range("I12").Value = 0.02758982345723265
$ws.Range("J12").Value = 0.02758982345723265
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 0.829094
$ws.Range("N12").Value = 2.487282
$ws.Range("O12").Value = 0.03374004171190829
$ws.Range("P12").Value = 0.03374004171190828
$ws.Range("Q12").Value = 0.31826845925
$ws.Range("R12").Value = 2.86441613325
$ws.Range("S12").Value = 0.0009308817942712155
$ws.Range("T12").Value = 0.0009308817942712153

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.383875
$ws.Range("H13").Value = 1.151625
$ws.Range("I13").Value = 0.02758982345723265
$ws.Range("J13").Value = 0.02758982345723265
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 9.203668
$ws.Range("N13").Value = 27.611004
$ws.Range("O13").Value = 0.3745439506528278
$ws.Range("P13").Value = 0.3745439506528276
$ws.Range("Q13").Value = 3.5330580535
$ws.Range("R13").Value = 31.7975224815
$ws.Range("S13").Value = 0.01033360147548598
$ws.Range("T13").Value = 0.01033360147548597

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 2.264687
$ws.Range("H14").Value = 6.794061
$ws.Range("I14").Value = 0.162767344880208
$ws.Range("J14").Value = 0.162767344880208
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 1.378475333333333
$ws.Range("N14").Value = 4.135426
$ws.Range("O14").Value = 0.05609715574531157
$ws.Range("P14").Value = 0.05609715574531156
$ws.Range("Q14").Value = 3.121815167220666
$ws.Range("R14").Value = 28.096336504986
$ws.Range("S14").Value = 0.00913078509599587
$ws.Range("T14").Value = 0.009130785095995868

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 2.264687
$ws.Range("H15").Value = 6.794061
$ws.Range("I15").Value = 0.162767344880208
$ws.Range("J15").Value = 0.162767344880208
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 13.16176133333333
$ws.Range("N15").Value = 39.485284
$ws.Range("O15").Value = 0.5356188518899525
$ws.Range("P15").Value = 0.5356188518899525
$ws.Range("Q15").Value = 29.80726978870267
$ws.Range("R15").Value = 268.265428098324
$ws.Range("S15").Value = 0.08718125838991295
$ws.Range("T15").Value = 0.08718125838991295

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 2.264687
$ws.Range("H16").Value = 6.794061
$ws.Range("I16").Value = 0.162767344880208
$ws.Range("J16").Value = 0.162767344880208
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 0.829094
$ws.Range("N16").Value = 2.487282
$ws.Range("O16").Value = 0.03374004171190829
$ws.Range("P16").Value = 0.03374004171190828
$ws.Range("Q16").Value = 1.877638403578
$ws.Range("R16").Value = 16.898745632202
$ws.Range("S16").Value = 0.00549177700559478
$ws.Range("T16").Value = 0.005491777005594779

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 2.264687
$ws.Range("H17").Value = 6.794061
$ws.Range("I17").Value = 0.162767344880208
$ws.Range("J17").Value = 0.162767344880208
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 9.203668
$ws.Range("N17").Value = 27.611004
$ws.Range("O17").Value = 0.3745439506528278
$ws.Range("P17").Value = 0.3745439506528276
$ws.Range("Q17").Value = 20.843427271916
$ws.Range("R17").Value = 187.590845447244
$ws.Range("S17").Value = 0.06096352438870442
$ws.Range("T17").Value = 0.0609635243887044
